# "Novo game para teste de funcionalidades baseado no GiseUp"
#
# 1) Slide 1: ungroup "Grupo 3" (the floating card decoration in the
#    top-right corner) into its four member shapes and shift the whole
#    cluster down-and-left a bit.
# 2) Slide 2: nudge the picture placeholder into its new spot.
# 3) Slide 2: clone the blue "card" rectangle from slide 1 onto slide 2
#    (recolored to a flat blue) as a new "Retângulo 4" shape.

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$s2 = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# 1) Ungroup "Grupo 3" on slide 1 and move the pieces to their new spot.
# ---------------------------------------------------------------------
$grupo3 = $s1.Shapes.Item(1)
$pieces = $grupo3.Ungroup()

foreach ($shp in $pieces) {
    switch ($shp.Id) {
        27 { $shp.Left = 756.643951507874;  $shp.Top = 120.45393760787402 }  # Retângulo 26
        28 { $shp.Left = 726.958188976378;  $shp.Top = 254.75897987795275 }  # Menos 27
        22 { $shp.Left = 795.2511291622047; $shp.Top = 73.48031496062993 }   # Semicírculos 21
        23 { $shp.Left = 795.4848938897637; $shp.Top = 363.14897637795275 } # Semicírculos 22
    }
}

# Keep a handle on the (still blip-filled) "Retângulo 26" shape — it is
# used below as the template for the new rectangle added to slide 2.
$rect26 = $null
foreach ($shp in $pieces) {
    if ($shp.Id -eq 27) { $rect26 = $shp }
}

# ---------------------------------------------------------------------
# 2) Slide 2: reposition the picture ("Imagem 1").
# ---------------------------------------------------------------------
$imagem1 = $s2.Shapes.Item(1)
$imagem1.Left = 536.5874939149606
$imagem1.Top = 91.34299212598425

# ---------------------------------------------------------------------
# 3) Slide 2: add the new "Retângulo 4" shape — a copy of "Retângulo 26"
#    (so it inherits the same style/txBody), recolored to a flat blue
#    solid fill instead of the picture fill.
# ---------------------------------------------------------------------
$rect26.Copy()
$pasted = $s2.Shapes.Paste()
$newRect = $pasted.Item(1)

$newRect.Name = "Retângulo 4"
$newRect.Fill.Solid()
$newRect.Fill.ForeColor.RGB = 0xF0B000
$newRect.Left = 54.49055298110236
$newRect.Top = 153.53504187007874
